$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A1: "NOM DU USE CASE : supprimerquestionner" -> "NOM DU USE CASE : supprimersondage" ---
# Rich text: run1 "NOM DU USE CASE" (bold, 12pt), run2 " : supprimerquestionner" (regular, 11pt)
$cellA1 = $ws.Range("A1")
$oldTail1 = " : supprimerquestionner"
$newTail1 = " : supprimersondage"
$cellA1.Characters(16, $oldTail1.Length).Text = $newTail1

$a1r1 = $cellA1.Characters(1, 15)
$a1r1.Font.Bold = $true
$a1r1.Font.Size = 12

$a1r2 = $cellA1.Characters(16, $newTail1.Length)
$a1r2.Font.Bold = $false
$a1r2.Font.Size = 11

# --- A2: "...BUT : supprimer un question questionnaire" -> "...BUT : supprimer un  sondage" ---
# Rich text: run1 "ACTEUR" (bold,12), run2 " : Admin, super admin\n\n" (regular,11),
#            run3 "BUT" (bold,12), run4 " : supprimer un question questionnaire" (regular,11)
$cellA2 = $ws.Range("A2")
$oldTail2 = " : supprimer un question questionnaire"
$newTail2 = " : supprimer un  sondage"
$cellA2.Characters(33, $oldTail2.Length).Text = $newTail2

$a2r1 = $cellA2.Characters(1, 6)
$a2r1.Font.Bold = $true
$a2r1.Font.Size = 12

$a2r2 = $cellA2.Characters(7, 23)
$a2r2.Font.Bold = $false
$a2r2.Font.Size = 11

$a2r3 = $cellA2.Characters(30, 3)
$a2r3.Font.Bold = $true
$a2r3.Font.Size = 12

$a2r4 = $cellA2.Characters(33, $newTail2.Length)
$a2r4.Font.Bold = $false
$a2r4.Font.Size = 11

# --- A6: "2) Quideance supprime le questionnaire" -> "2) Quideance supprime le sondage" ---
$ws.Range("A6").Value = "2) Quideance supprime le sondage"
